$members = $excel | Get-Member
Write-Host $members
